# Applies the commit diff:
#   1. "Вариант 2" -> "Вариант 1"
#   2. Moves the "_GoBack" bookmark from right after "LinkedList.cpp"
#      to right after the (now) "Вариант 1" run.
#
# Note: this runtime's Bookmarks.Add() has a quirk where a collapsed
# range sitting exactly one character before a paragraph mark gets
# mis-resolved (it snaps back to the very start of the document).
# To work around it we temporarily pad the target spot with two
# placeholder characters, drop the bookmark safely between them (which
# is no longer adjacent to the paragraph mark), and then remove the
# padding again - the collapsed bookmark stays anchored to the correct
# text position once the padding disappears.

$d = $word.ActiveDocument

# --- 1. Rename "Вариант 2" -> "Вариант 1" ---------------------------------
$rng = $d.Content
$rng.Find.Execute("Вариант 2", $true, $false, $false, $false, $false, $true, `
                   1, $false, "Вариант 1", 2) | Out-Null

# --- 2. Drop the old "_GoBack" bookmark (after "LinkedList.cpp") ---------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- 3. Re-create "_GoBack" right after the "Вариант 1" run ---------------
$target = $d.Content
$target.Find.Execute("Вариант 1", $true, $false, $false, $false, $false, `
                      $true, 1, $false, "", 0) | Out-Null
$target.Collapse(0)

# Pad with two throw-away characters so the insertion point used for the
# bookmark is no longer flush against the paragraph mark.
$target.InsertAfter("XY")
$target.Collapse(0)
$target.MoveStart(1, -1) | Out-Null
$target.Collapse(1)

$d.Bookmarks.Add("_GoBack", $target)
$bm = $d.Bookmarks("_GoBack")

# Remove the padding: the character right after the bookmark ("Y"), then
# the character right before it ("X"). The collapsed bookmark stays put.
$after = $d.Range($bm.Start, $bm.Start + 1)
$after.Text = ""

$bm = $d.Bookmarks("_GoBack")
$before = $d.Range($bm.Start - 1, $bm.Start)
$before.Text = ""
